$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KN Curves")

$headers = @("Displacement","5°","10°","12°","15°","20°","25°","30°","35°","40°","50°","60°","75°","90°")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$data = @(
    @(10000,3.5,7,8.2,9.9,11.45,11.85,11.4,11,10.4,9.4,8.7,8.5,8.6),
    @(20000,2.15,4.2,5,6.2,8,9.6,10.75,11.45,11.85,11.9,11.45,10.5,9.8),
    @(30000,1.6,3.15,3.8,4.7,6.2,7.55,8.8,9.9,10.85,11.75,11.9,11.45,10.8),
    @(40000,1.4,2.75,3.3,4.05,5.35,6.55,7.65,8.6,9.55,10.9,11.55,11.85,11.45),
    @(50000,1.3,2.55,3.05,3.75,4.95,5.95,6.95,7.85,8.7,10.05,10.95,11.6,11.8),
    @(60000,1.25,2.5,2.95,3.6,4.7,5.6,6.5,7.35,8.15,9.45,10.4,11.15,11.5),
    @(70000,1.2,2.45,2.9,3.45,4.5,5.35,6.15,6.95,7.7,8.95,9.85,10.7,11.1),
    @(80000,1.2,2.4,2.85,3.4,4.35,5.1,5.9,6.65,7.35,8.5,9.4,10.25,10.7),
    @(90000,1.2,2.4,2.85,3.4,4.25,5,5.7,6.4,7.05,8.15,9,9.8,10.25),
    @(100000,1.2,2.35,2.8,3.35,4.15,4.9,5.55,6.2,6.85,7.85,8.65,9.5,9.95),
    @(110000,1.2,2.3,2.75,3.25,4.1,4.8,5.4,6,6.65,7.6,8.4,9.2,9.65)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowData = $data[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

# Clear any leftover rows/columns beyond the new data range (previously up to L51)
$ws.Range("A13:L51").Clear()
